$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 2 (the THIAGO row), pushing existing data down.
$ws.Rows.Item(2).Resize(4).Insert()

# New account rows to add, in order, right above THIAGO.
$newRows = @(
    @("004352384", "BRASFORT", 109482.35),
    @("005870700", "ALOISIO", 100000),
    @("002973105", "DARLAN", 50859.58),
    @("004216504", "WANDER", 41448.73)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i

    # Write the account number as a formula that evaluates to text, then
    # paste back as a value-only copy. This keeps the leading zeros (the
    # field is textual, like "004352384") without leaving a residual
    # formula behind and without altering the cell's number format/style
    # (plain `.Value = "004352384"` would be auto-coerced to the number
    # 4352384, losing the leading zeros).
    $ws.Cells.Item($r, 1).Formula = '="' + $newRows[$i][0] + '"'
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# Remove the GISELA row, which now sits right after THIAGO (row 6 => row 7).
$ws.Rows.Item(7).Delete()
